$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-case rows 44-49 appended to the "SendQuote" block ---
# Shared strings must be appended in a specific first-use order so the
# resulting sharedStrings.xml unique-string order matches the target
# workbook (indices 71..76). That order is NOT simple row order, so we
# deliberately touch the cells in the sequence that yields:
#   71: 102_AutomobileInsurance_006_SendQuote_002_EnterValuesInWrongFormat
#   72: 102_AutomobileInsurance_006_SendQuote_002_EnterValuesInWrongFormat Part 2
#   73: Send Quote Page check error hint formatting Part 2
#   74: Send Quote Page check error hint formatting
#   75: Send Quote page check for open mandatory field
#   76: Send Quote Page check for hints regarding mandatory fields

$txt71 = "102_AutomobileInsurance_006_SendQuote_002_EnterValuesInWrongFormat"
$txt72 = "102_AutomobileInsurance_006_SendQuote_002_EnterValuesInWrongFormat Part 2"
$txt73 = "Send Quote Page check error hint formatting Part 2"
$txt74 = "Send Quote Page check error hint formatting"
$txt75 = "Send Quote page check for open mandatory field"
$txt76 = "Send Quote Page check for hints regarding mandatory fields"

$CHK = "<CHK>"
$SET = "<SET>"
$NOP = "<NOP>"

# Row 46 first (introduces txt71)
$ws.Cells.Item(46,1).Value = $txt71
$ws.Cells.Item(46,1).NumberFormat = "@"
$ws.Cells.Item(46,2).Value = $SET
$ws.Cells.Item(46,7).Value = $txt71
$ws.Cells.Item(46,7).NumberFormat = "@"
$ws.Cells.Item(46,8).Value = $NOP

# Row 48 (introduces txt72)
$ws.Cells.Item(48,1).Value = $txt72
$ws.Cells.Item(48,1).NumberFormat = "@"
$ws.Cells.Item(48,2).Value = $SET
$ws.Cells.Item(48,7).Value = $txt72
$ws.Cells.Item(48,7).NumberFormat = "@"
$ws.Cells.Item(48,8).Value = $NOP

# Row 49 (introduces txt73)
$ws.Cells.Item(49,1).Value = $txt73
$ws.Cells.Item(49,1).NumberFormat = "@"
$ws.Cells.Item(49,2).Value = $CHK
$ws.Cells.Item(49,7).Value = $txt73
$ws.Cells.Item(49,7).NumberFormat = "@"
$ws.Cells.Item(49,8).Value = $NOP

# Row 47 (introduces txt74)
$ws.Cells.Item(47,1).Value = $txt74
$ws.Cells.Item(47,1).NumberFormat = "@"
$ws.Cells.Item(47,2).Value = $CHK
$ws.Cells.Item(47,7).Value = $txt74
$ws.Cells.Item(47,7).NumberFormat = "@"
$ws.Cells.Item(47,8).Value = $NOP

# Row 44 (introduces txt75)
$ws.Cells.Item(44,1).Value = $txt75
$ws.Cells.Item(44,2).Value = $CHK
$ws.Cells.Item(44,7).Value = $txt75
$ws.Cells.Item(44,8).Value = $NOP

# Row 45 (introduces txt76)
$ws.Cells.Item(45,1).Value = $txt76
$ws.Cells.Item(45,2).Value = $CHK
$ws.Cells.Item(45,7).Value = $txt76
$ws.Cells.Item(45,8).Value = $NOP

# --- Move the "Grafik 2" picture below the new rows ---
$shp = $ws.Shapes.Item(1)
$shp.Top = 737.4000787401575

# --- Update the active selection to match the end-state view ---
$ws.Range("G48").Select()
